$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 68628220
$ws.Range("I43").Value = 100000400
$ws.Range("J43").Value = 23810838
$ws.Range("K43").Value = 100000400
$ws.Range("L43").Value = 23810838
$ws.Range("M43").Value = -100000331
$ws.Range("N43").Value = -23810976
$ws.Range("H92").Value = 1615.3158
$ws.Range("I92").Value = 1761.2307
$ws.Range("K92").Value = 1761.2307
$ws.Range("M92").Value = -513.2307000000001
$ws.Range("H112").Value = 1285.75
$ws.Range("J112").Value = 1327.1052
$ws.Range("L112").Value = 3981.3156
$ws.Range("N112").Value = -6197.3156
$ws.Range("H132").Value = 2422565.2
$ws.Range("I132").Value = 2696527.2
$ws.Range("K132").Value = 8089581.600000001
$ws.Range("M132").Value = -8087051.600000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3425.4285
$ws.Range("I32").Value = 2578.5647
$ws.Range("J32").Value = 8962.615
$ws.Range("K32").Value = 2578.5647
$ws.Range("L32").Value = 8962.615
$ws.Range("M32").Value = -2291.5647
$ws.Range("N32").Value = -9536.615
$ws.Range("H45").Value = 1104.9131
$ws.Range("I45").Value = 1077.1818
$ws.Range("K45").Value = 1077.1818
$ws.Range("M45").Value = -700.1818000000001
$ws.Range("H63").Value = 2501875
$ws.Range("I63").Value = 10000000
$ws.Range("J63").Value = 2500
$ws.Range("K63").Value = 10000000
$ws.Range("L63").Value = 2500
$ws.Range("M63").Value = -9999314
$ws.Range("N63").Value = -3872
$ws.Range("H66").Value = 2501875
$ws.Range("I66").Value = 10000000
$ws.Range("J66").Value = 2500
$ws.Range("K66").Value = 50000000
$ws.Range("L66").Value = 12500
$ws.Range("M66").Value = -49996568
$ws.Range("N66").Value = -19364
$ws.Range("H111").Value = 30000
$ws.Range("J111").Value = 30000
$ws.Range("L111").Value = 30000
$ws.Range("N111").Value = -38180
$ws.Range("H122").Value = 2297.4
$ws.Range("I122").Value = 2297.4
$ws.Range("K122").Value = 6892.200000000001
$ws.Range("M122").Value = -4442.200000000001
$ws.Range("H132").Value = 3689.9614
$ws.Range("I132").Value = 3667.8445
$ws.Range("J132").Value = 3832.1428
$ws.Range("K132").Value = 11003.5335
$ws.Range("L132").Value = 11496.4284
$ws.Range("M132").Value = -8473.533500000001
$ws.Range("N132").Value = -16556.4284

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H128").Value = 3000
$ws.Range("I128").Value = 3000
$ws.Range("K128").Value = 9000
$ws.Range("M128").Value = -6510
$ws.Range("H134").Value = 36299.766
$ws.Range("I134").Value = 94404.63
$ws.Range("J134").Value = 2660.1052
$ws.Range("K134").Value = 283213.89
$ws.Range("L134").Value = 7980.3156
$ws.Range("M134").Value = -280678.89
$ws.Range("N134").Value = -13050.3156

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5053608.5
$ws.Range("I31").Value = 2334.9524
$ws.Range("K31").Value = 2334.9524
$ws.Range("M31").Value = -2039.9524
$ws.Range("H34").Value = 5053608.5
$ws.Range("I34").Value = 2334.9524
$ws.Range("K34").Value = 2334.9524
$ws.Range("M34").Value = -2132.9524
$ws.Range("H86").Value = 2455.75
$ws.Range("I86").Value = 1174.625
$ws.Range("J86").Value = 3736.875
$ws.Range("K86").Value = 1174.625
$ws.Range("L86").Value = 3736.875
$ws.Range("M86").Value = -51.625
$ws.Range("N86").Value = -5982.875
$ws.Range("H89").Value = 2455.75
$ws.Range("I89").Value = 1174.625
$ws.Range("J89").Value = 3736.875
$ws.Range("K89").Value = 5873.125
$ws.Range("L89").Value = 18684.375
$ws.Range("M89").Value = -257.125
$ws.Range("N89").Value = -29916.375
$ws.Range("H132").Value = 2402.5945
$ws.Range("I132").Value = 1608.5217
$ws.Range("J132").Value = 3707.1428
$ws.Range("K132").Value = 4825.5651
$ws.Range("L132").Value = 11121.4284
$ws.Range("M132").Value = -2295.5651
$ws.Range("N132").Value = -16181.4284
$ws.Range("H134").Value = 879.36957
$ws.Range("I134").Value = 803.85364
$ws.Range("J134").Value = 1498.6
$ws.Range("K134").Value = 2411.56092
$ws.Range("L134").Value = 4495.799999999999
$ws.Range("M134").Value = 123.4390800000001
$ws.Range("N134").Value = -9565.799999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1347.9
$ws.Range("I4").Value = 1075
$ws.Range("J4").Value = 1416.125
$ws.Range("K4").Value = 3225
$ws.Range("L4").Value = 4248.375
$ws.Range("M4").Value = -3113
$ws.Range("N4").Value = -4472.375
$ws.Range("H60").Value = 513.3889
$ws.Range("I60").Value = 383.8889
$ws.Range("J60").Value = 642.8889
$ws.Range("K60").Value = 1151.6667
$ws.Range("L60").Value = 1928.6667
$ws.Range("M60").Value = -900.6667
$ws.Range("N60").Value = -2430.6667
$ws.Range("H92").Value = 944.4167
$ws.Range("I92").Value = 1237.5
$ws.Range("J92").Value = 797.875
$ws.Range("K92").Value = 3712.5
$ws.Range("L92").Value = 2393.625
$ws.Range("M92").Value = -2464.5
$ws.Range("N92").Value = -4889.625
$ws.Range("H107").Value = 419.66666
$ws.Range("I107").Value = 772.5
$ws.Range("J107").Value = 291.36365
$ws.Range("K107").Value = 2317.5
$ws.Range("L107").Value = 874.09095
$ws.Range("M107").Value = -397.5
$ws.Range("N107").Value = -4714.09095
$ws.Range("H114").Value = 3458.3076
$ws.Range("I114").Value = 2757
$ws.Range("J114").Value = 3770
$ws.Range("K114").Value = 8271
$ws.Range("L114").Value = 11310
$ws.Range("M114").Value = -5017
$ws.Range("N114").Value = -17818

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 31250782
$ws.Range("I113").Value = 50000436
$ws.Range("J113").Value = 1360
$ws.Range("K113").Value = 50000436
$ws.Range("L113").Value = 1360
$ws.Range("M113").Value = -49998266
$ws.Range("N113").Value = -5700
$ws.Range("H132").Value = 44396.668
$ws.Range("I132").Value = 54605.684
$ws.Range("J132").Value = 5602.4
$ws.Range("K132").Value = 163817.052
$ws.Range("L132").Value = 16807.2
$ws.Range("M132").Value = -161287.052
$ws.Range("N132").Value = -21867.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H56").Value = 35000
$ws.Range("I56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("M56").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("M55").ClearContents()
$ws.Range("H122").Value = 53000.55
$ws.Range("I122").Value = 61773.65
$ws.Range("K122").Value = 185320.95
$ws.Range("M122").Value = -182870.95
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H126").Value = 7958.5884
$ws.Range("I126").Value = 9646
$ws.Range("J126").Value = 2474.5
$ws.Range("K126").Value = 28938
$ws.Range("L126").Value = 7423.5
$ws.Range("M126").Value = -26468
$ws.Range("N126").Value = -12363.5

